# Weekly update: insert a new price row for the current week (row 216),
# pushing the existing history down by one row (216->217, ..., 222->223).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 216, shifting rows 216:222
# down to 217:223 (and the sheet dimension grows from R222 to R223).
$ws.Rows.Item(216).Insert()

# Populate the newly inserted row 216 with this week's data.
$ws.Cells.Item(216, 1).Value = 7
$ws.Cells.Item(216, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(216, 3).Value = "Ñuble"
$ws.Cells.Item(216, 4).Value = 45041
$ws.Cells.Item(216, 5).Value = 16
$ws.Cells.Item(216, 6).Value = 100112040
$ws.Cells.Item(216, 7).Value = "Cilantro"
$ws.Cells.Item(216, 8).Value = "Sin especificar"
$ws.Cells.Item(216, 9).Value = "Primera"
$ws.Cells.Item(216, 10).Value = 100
$ws.Cells.Item(216, 11).Value = 1500
$ws.Cells.Item(216, 12).Value = 1500
$ws.Cells.Item(216, 13).Value = 1500
$ws.Cells.Item(216, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(216, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(216, 16).Value = 1500
$ws.Cells.Item(216, 17).Value = 1
$ws.Cells.Item(216, 18).Value = "Hortaliza"
